# Initial integration of label gateway
#
# - Rename "ValueSwapTable" -> "RegionFixTable", give it a small 2-row
#   lookup table (No / Region / Fix), make it the active sheet.
# - Add a brand-new "ValueFixTable" sheet after it, with a 4-row lookup
#   table (No / Value / Fix) plus a trailing numbered-but-otherwise-empty
#   column A down to 34.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. RegionFixTable (rename of the old, empty ValueSwapTable)
# ---------------------------------------------------------------------
$regionFix = $wb.Worksheets.Item("ValueSwapTable")
$regionFix.Name = "RegionFixTable"

$regionFix.Range("A1").Value = "No"
$regionFix.Range("B1").Value = "Region"
$regionFix.Range("C1").Value = "Fix"
$regionFix.Range("A1:C1").Font.Bold = $true

$regionFix.Range("A2").Value = 1
$regionFix.Range("B2").Value = "World"
$regionFix.Range("C2").Value = "WLD"

$regionFix.Columns.Item(1).ColumnWidth = 2.6667
$regionFix.Columns.Item(2).ColumnWidth = 15.5567
$regionFix.Columns.Item(3).ColumnWidth = 18.4167

# ---------------------------------------------------------------------
# 2. ValueFixTable (new sheet, placed right after RegionFixTable)
# ---------------------------------------------------------------------
$valueFix = $wb.Worksheets.Add($null, $regionFix)
$valueFix.Name = "ValueFixTable"

$valueFix.Range("A1").Value = "No"
$valueFix.Range("B1").Value = "Value"
$valueFix.Range("C1").Value = "Fix"
$valueFix.Range("A1:C1").Font.Bold = $true

# Leading "'" forces these to be stored as literal text instead of being
# auto-coerced into an Excel error value / other special type.
$fixLabels = @("'#div/0!", "n/a", "na", "nan")
for ($i = 0; $i -lt $fixLabels.Length; $i++) {
    $r = $i + 2
    $valueFix.Cells.Item($r, 1).Value = $i + 1
    $valueFix.Cells.Item($r, 2).Value = $fixLabels[$i]
    $valueFix.Cells.Item($r, 3).Value = 0
}

# Rows 6..35 only carry the running index in column A.
for ($r = 6; $r -le 35; $r++) {
    $valueFix.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------
# 3. Active tab moves to the (renamed) RegionFixTable sheet.
# ---------------------------------------------------------------------
[void]$regionFix.Activate()
[void]$regionFix.Range("A3").Select()
